$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): F3 2975 -> 2982, F5 73 -> 74
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2982
$ws1.Range("F5").Value = 74

# Sheet "全部类型" (All types): F7 2975 -> 2982, F10 73 -> 74
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2982
$ws4.Range("F10").Value = 74
